# Applies reordering of dict keys within the text representations stored
# in column A of the "Subgroups" worksheet, as described by the commit
# "added my own version of random walks".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$updates = @{
    10 = "{'Hobby': np.int64(1), 'Student': np.int64(1)}"
    11 = "{'Hobby': np.int64(1), 'SexualOrientation': np.int64(1)}"
    13 = "{'Dependents': np.int64(2), 'Hobby': np.int64(1)}"
    17 = "{'HDI': np.int64(1), 'Student': np.int64(1)}"
    18 = "{'UndergradMajor': np.int64(2), 'SexualOrientation': np.int64(1)}"
    21 = "{'Dependents': np.int64(2), 'SexualOrientation': np.int64(1)}"
    22 = "{'HDI': np.int64(1), 'SexualOrientation': np.int64(1)}"
    24 = "{'Hobby': np.int64(1), 'Student': np.int64(1), 'SexualOrientation': np.int64(1)}"
    25 = "{'RaceEthnicity': np.int64(1), 'Hobby': np.int64(1), 'SexualOrientation': np.int64(1)}"
    26 = "{'Hobby': np.int64(1), 'SexualOrientation': np.int64(1), 'HDI': np.int64(1)}"
    28 = "{'HDI': np.int64(1), 'Student': np.int64(1), 'SexualOrientation': np.int64(1)}"
    29 = "{'RaceEthnicity': np.int64(1), 'HDI': np.int64(1), 'Student': np.int64(1)}"
    30 = "{'RaceEthnicity': np.int64(1), 'HDI': np.int64(1), 'SexualOrientation': np.int64(1)}"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
